# Differentiate dermis (exodermis/endodermis) conductances according to apex distance.
# scenarios_as_columns: columns I:N (scenarios 3-8) get their own values/formulas,
# decoupled from the shared formulas that used to mirror column H across the board.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios_as_columns")
$ws.Activate()

# Rows 106-109: give columns I:N their own formulas (no longer copies of column H)
$ws.Range("I106:N106").Formula = "=10/0.1"
$ws.Range("I107:N107").Formula = "=150/0.1"
$ws.Range("I108:N108").Formula = "=20/0.1"
$ws.Range("I109:N109").Formula = "=400/0.1"

# Row 209: relative conductance at the apex -- lowered for the I:N scenarios
$ws.Range("I209:N209").Value = 1.36

# Row 210: associated distance parameter -- raised (less negative) for the I:N scenarios
# (row 211 mirrors row 210 through its own formula and recalculates automatically)
$ws.Range("I210:N210").Value = -30000

# Restore the view to where the author was last working (row ~108, column I)
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I108").Select()
$excel.ActiveWindow.ScrollRow = 102
$excel.ActiveWindow.ScrollColumn = 5
